$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank data rows right above the current row 93, shifting all
# the existing data (old rows 93-150) down to rows 95-152.
$ws.Rows.Item(93).EntireRow.Insert()
$ws.Rows.Item(93).EntireRow.Insert()

# New row 93: Early Majestic / Primera, Provincia de Limari, $/bandeja 10 kilos granel
$ws.Cells.Item(93, 1).Value = 10
$ws.Cells.Item(93, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(93, 3).Value = "La Araucanía"
$ws.Cells.Item(93, 4).Value = 44529
$ws.Cells.Item(93, 5).Value = 9
$ws.Cells.Item(93, 6).Value = "Fruta"
$ws.Cells.Item(93, 7).Value = 100103
$ws.Cells.Item(93, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(93, 9).Value = 100103004
$ws.Cells.Item(93, 10).Value = "Durazno"
$ws.Cells.Item(93, 11).Value = "Early Majestic"
$ws.Cells.Item(93, 12).Value = "Primera"
$ws.Cells.Item(93, 13).Value = 210
$ws.Cells.Item(93, 14).Value = 13000
$ws.Cells.Item(93, 15).Value = 13000
$ws.Cells.Item(93, 16).Value = 13000
$ws.Cells.Item(93, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(93, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(93, 19).Value = 1300
$ws.Cells.Item(93, 20).Value = 10

# New row 94: Early Majestic / Primera, Provincia de Limari, $/bandeja 15 kilos granel
$ws.Cells.Item(94, 1).Value = 10
$ws.Cells.Item(94, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(94, 3).Value = "La Araucanía"
$ws.Cells.Item(94, 4).Value = 44529
$ws.Cells.Item(94, 5).Value = 9
$ws.Cells.Item(94, 6).Value = "Fruta"
$ws.Cells.Item(94, 7).Value = 100103
$ws.Cells.Item(94, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(94, 9).Value = 100103004
$ws.Cells.Item(94, 10).Value = "Durazno"
$ws.Cells.Item(94, 11).Value = "Early Majestic"
$ws.Cells.Item(94, 12).Value = "Primera"
$ws.Cells.Item(94, 13).Value = 35
$ws.Cells.Item(94, 14).Value = 23000
$ws.Cells.Item(94, 15).Value = 23000
$ws.Cells.Item(94, 16).Value = 23000
$ws.Cells.Item(94, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(94, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(94, 19).Value = 1533
$ws.Cells.Item(94, 20).Value = 15
